# Update the cryptos worksheet with refreshed price / 1h-volume figures and
# reorder a few rows (Decentraland/EnergySwap and BabyDogeCoin/PancakeSwap
# swapped places), matching the latest GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values look like plain numbers (e.g. "1.003",
# "334.83"). The source data stores these as text (inlineStr) so that
# formatting such as trailing zeros is preserved exactly. Force the
# NumberFormat to Text ("@") on just those cells before assigning the
# value, otherwise Excel would silently convert them into numeric values
# and drop significant trailing zeros / change the cell type.
$textForceCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '30.452.97'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '2.105.64'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '334.83'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.5245'
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("D8").Value = '0.4605'
$ws.Range("E8").Value = '  +6.18%  '
$ws.Range("D9").Value = '52.43'
$ws.Range("E9").Value = '  +12.04%  '
$ws.Range("D10").Value = '0.08959'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").Value = '1.178'
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").Value = '24.42'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '2.094.60'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '6.791'
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("D15").Value = '7.895'
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '96.32'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '0.00001132'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").Value = '0.06628'
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("D20").Value = '19.28'
$ws.Range("E20").Value = '  +1.85%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '6.282'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '30.513.27'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("D24").Value = '12.33'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").Value = '2.366'
$ws.Range("E25").Value = '  +3.49%  '
$ws.Range("D26").Value = '2.337.23'
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("D27").Value = '22.29'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = '2.567'
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").Value = '163.54'
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("D30").Value = '132.63'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '1.196'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '0.1073'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").Value = '1.698'
$ws.Range("E33").Value = '  +10.39%  '
$ws.Range("D34").Value = '6.155'
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").Value = '3.921'
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("E36").Value = '  +8.93%  '
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("D38").Value = '0.06822'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = '5.550'
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").Value = '12.82'
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("D41").Value = '0.2287'
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").Value = '0.6893'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("D43").Value = '1.248'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '2.337'
$ws.Range("E44").Value = '  +5.93%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6386'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '13.95'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.00000000360'
$ws.Range("E48").Value = '  +24.79%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '3.656'
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").Value = '1.214'
$ws.Range("E51").Value = '  +1.60%  '
